$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.921492695808411
$ws.Range("B1").Value = 2.954934120178223
$ws.Range("C1").Value = 2.359329462051392
$ws.Range("D1").Value = 2.230969429016113
$ws.Range("E1").Value = 1.930059194564819
